$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 1224
$ws.Range("E2").Value = -14
$ws.Range("F2").Value = -14
$ws.Range("G2").Value = -57
$ws.Range("H2").Value = -69
$ws.Range("I2").Value = -69
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 1175
$ws.Range("L2").Value = 667
$ws.Range("M2").Value = 509
$ws.Range("N2").Value = 508
$ws.Range("O2").Value = 1
$ws.Range("P2").Value = 159
$ws.Range("Q2").Value = -1
$ws.Range("R2").Value = -6
$ws.Range("S2").Value = 4
$ws.Range("T2").Value = 24
$ws.Range("U2").Value = -24
$ws.Range("V2").Value = 354
$ws.Range("W2").Value = -1.16
$ws.Range("X2").Value = -5.65
$ws.Range("Y2").Value = -12.71
$ws.Range("Z2").Value = -5.89
$ws.Range("AA2").Value = 131.07
$ws.Range("AB2").Value = 229.07
$ws.Range("AC2").Value = -377
$ws.Range("AD2").Value = -5.63
$ws.Range("AE2").Value = 2935
$ws.Range("AF2").Value = 0.72
$ws.Range("AG2").Value = 9
$ws.Range("AH2").Value = 0.41
$ws.Range("AI2").Value = -2.13
$ws.Range("AJ2").Value = 18400520
$ws.Range("D3").Value = 1176
$ws.Range("E3").Value = 21
$ws.Range("F3").Value = 21
$ws.Range("G3").Value = -17
$ws.Range("H3").Value = -59
$ws.Range("I3").Value = -59
$ws.Range("J3").Value = 1
$ws.Range("K3").Value = 1366
$ws.Range("L3").Value = 682
$ws.Range("M3").Value = 685
$ws.Range("N3").Value = 683
$ws.Range("O3").Value = 2
$ws.Range("P3").Value = 166
$ws.Range("Q3").Value = -51
$ws.Range("R3").Value = -17
$ws.Range("S3").Value = 67
$ws.Range("T3").Value = 12
$ws.Range("U3").Value = -63
$ws.Range("V3").Value = 404
$ws.Range("W3").Value = 1.78
$ws.Range("X3").Value = -5
$ws.Range("Y3").Value = -9.98
$ws.Range("Z3").Value = -4.63
$ws.Range("AA3").Value = 99.59999999999999
$ws.Range("AB3").Value = 181.73
$ws.Range("AC3").Value = -323
$ws.Range("AD3").Value = -10.64
$ws.Range("AE3").Value = 3834
$ws.Range("AF3").Value = 0.9
$ws.Range("AG3").Value = 0
$ws.Range("AH3").Value = 0
$ws.Range("AI3").Value = 0
$ws.Range("AJ3").Value = 18400520
$ws.Range("D4").Value = 1207
$ws.Range("E4").Value = 40
$ws.Range("F4").Value = 40
$ws.Range("G4").Value = 18
$ws.Range("H4").Value = 8
$ws.Range("I4").Value = 8
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 1295
$ws.Range("L4").Value = 619
$ws.Range("M4").Value = 675
$ws.Range("N4").Value = 675
$ws.Range("O4").Value = 0
$ws.Range("P4").Value = 166
$ws.Range("Q4").Value = 161
$ws.Range("R4").Value = -49
$ws.Range("S4").Value = -53
$ws.Range("T4").Value = 19
$ws.Range("U4").Value = 142
$ws.Range("V4").Value = 331
$ws.Range("W4").Value = 3.29
$ws.Range("X4").Value = 0.66
$ws.Range("Y4").Value = 1.14
$ws.Range("Z4").Value = 0.6
$ws.Range("AA4").Value = 91.73
$ws.Range("AB4").Value = 183.32
$ws.Range("AC4").Value = 42
$ws.Range("AD4").Value = 103.76
$ws.Range("AE4").Value = 3840
$ws.Range("AF4").Value = 1.14
$ws.Range("AG4").Value = 9
$ws.Range("AH4").Value = 0.21
$ws.Range("AI4").Value = 20.42
$ws.Range("AJ4").Value = 18400520
$ws.Range("D5").Value = 1233
$ws.Range("E5").Value = 26
$ws.Range("F5").Value = 26
$ws.Range("G5").Value = 21
$ws.Range("H5").Value = 10
$ws.Range("I5").Value = 11
$ws.Range("J5").Value = -1
$ws.Range("K5").Value = 1316
$ws.Range("L5").Value = 639
$ws.Range("M5").Value = 678
$ws.Range("N5").Value = 677
$ws.Range("O5").Value = 0
$ws.Range("P5").Value = 174
$ws.Range("Q5").Value = 8
$ws.Range("R5").Value = -55
$ws.Range("S5").Value = 9
$ws.Range("T5").Value = 26
$ws.Range("U5").Value = -18
$ws.Range("V5").Value = 345
$ws.Range("W5").Value = 2.08
$ws.Range("X5").Value = 0.84
$ws.Range("Y5").Value = 1.62
$ws.Range("Z5").Value = 0.8
$ws.Range("AA5").Value = 94.23
$ws.Range("AB5").Value = 173.93
$ws.Range("AC5").Value = 59
$ws.Range("AD5").Value = 87.2
$ws.Range("AE5").Value = 3867
$ws.Range("AF5").Value = 1.34
$ws.Range("AG5").Value = 9
$ws.Range("AH5").Value = 0.18
$ws.Range("AI5").Value = 15.12
$ws.Range("AJ5").Value = 18400520
$ws.Range("D6").Value = 1077
$ws.Range("E6").Value = 33
$ws.Range("F6").Value = 33
$ws.Range("G6").Value = 28
$ws.Range("H6").Value = 22
$ws.Range("I6").Value = 22
$ws.Range("K6").Value = 1262
$ws.Range("L6").Value = 575
$ws.Range("M6").Value = 687
$ws.Range("N6").Value = 687
$ws.Range("P6").Value = 181
$ws.Range("Q6").Value = 132
$ws.Range("R6").Value = -64
$ws.Range("S6").Value = -57
$ws.Range("T6").Value = 30
$ws.Range("U6").Value = 102
$ws.Range("V6").Value = 290
$ws.Range("W6").Value = 3.05
$ws.Range("X6").Value = 2.03
$ws.Range("Y6").Value = 3.21
$ws.Range("Z6").Value = 1.69
$ws.Range("AA6").Value = 83.75
$ws.Range("AB6").Value = 170.98
$ws.Range("AC6").Value = 119
$ws.Range("AD6").Value = 35.83
$ws.Range("AE6").Value = 3921
$ws.Range("AG6").Value = 29
$ws.Range("AH6").Value = 0.6899999999999999
$ws.Range("AI6").Value = 23.49
$ws.Range("AJ6").Value = 18400520

# Clear out rows 7-9 (D through AJ), leaving only columns A-C populated
$ws.Range("D7:AJ9").ClearContents()
